$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13 ---
# (splits the old "Programa resumido:" row so that "Docentes responsaveis:"
#  at row 12 gets its own data row 13 for the professor name, and pushes
#  everything below down by one row)
$ws.Rows.Item(13).Insert()

# The insert carries a stray styled-but-empty cell into A13; clear it so the
# row has only the B/C value cells, matching the target layout.
$ws.Cells.Item(13, 1).Clear()

# --- Fix cell contents that changed / moved ---

# Row 10: Objetivos (PT) now holds the real Portuguese objectives text
# instead of the professor name that had been there by mistake.
$objPt = 'Compreender a Pesquisa Operacional como ciência aplicada à Engenharia de Produção. Proporcionar conhecimento dos problemas típicos encontrados em Engenharia de Produção. Analisar, Modelar e solucionar os problemas por meio da Pesquisa Operacional.'
$ws.Cells.Item(10, 2).Value = $objPt
$ws.Cells.Item(10, 3).Value = $objPt

# Row 13 (new): Docentes responsaveis value
$docente = '5840917 - Fabricio Maciel Gomes'
$ws.Cells.Item(13, 2).Value = $docente
$ws.Cells.Item(13, 3).Value = $docente
$ws.Cells.Item(13, 2).Font.Bold = $false
$ws.Cells.Item(13, 2).WrapText = $true

# Row 14: Programa resumido (PT short syllabus) replaces placeholder "Semestral"
$shortSylPt = 'Introdução a Pesquisa Operacional, Programação Linear, Método Simplex, Introdução aos Grafos e à Otimização em Rede, Estudo de Casos em Programação Linear, Introdução a Teoria das Filas,'
$ws.Cells.Item(14, 2).Value = $shortSylPt
$ws.Cells.Item(14, 3).Value = $shortSylPt

# Row 16: Programa (PT full syllabus) replaces placeholder date value
$fullSylPt = '1. Introdução a Pesquisa Operacional1.1. Conceitos de Pesquisa Operacional;1.2. Modelagem;1.3. Estrutura dos Modelos Matemáticos;1.4. Técnicas matemáticas em Pesquisa Operacional;1.2. Fases de Um Estudo em Pesquisa Operacional2. Programação Linear2.1. Definição2.2. Formulação de Modelos2.3. Resolução Gráfica;3. Método Simplex3.1. Desenvolvimento do Método Simplex;3.2. Procedimento do Método Simplex;4. Introdução aos Grafos e à Otimização em Rede4.1. Conceitos Básicos em Teoria dos Grafos4.2. Problemas de Fluxo Máximo;4.3. Problemas de Caminho Mínimo5. Estudo de Casos em Programação Linear5.1. Modelo de Transporte Simples5.2. Modelo da Designação.6. Introdução a Teoria das Filas6.1. Conceitos da Teoria das Filas6.2. Modelos Markovianos'
$ws.Cells.Item(16, 2).Value = $fullSylPt
$ws.Cells.Item(16, 3).Value = $fullSylPt

# Row 19: Metodo now holds the grading-method text instead of the professor name
$method = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Cells.Item(19, 2).Value = $method
$ws.Cells.Item(19, 3).Value = $method

# Row 20: Criterio now holds the passing-grade criterion text
$criterio = 'NF≥ 5,0.'
$ws.Cells.Item(20, 2).Value = $criterio
$ws.Cells.Item(20, 3).Value = $criterio

# Row 21: Norma de recuperacao now holds the make-up-exam rule text
$norma = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Cells.Item(21, 2).Value = $norma
$ws.Cells.Item(21, 3).Value = $norma

# Row 22: Bibliografia now holds the actual reading list
$biblio = "1. HILLIER, F.S., LIEBERMAN, G.J., Introdução à Pesquisa Operacional, 8ªed., Editora McGraw-Hill, 2006.`n2. LACHTERMACHER, G., Pesquisa Operacional na Tomada de Decisão (modelagem em Excel), 4ª ed., Editora Campus, 2009.`n3. ANDERSON, D.R., SWEENEY, D.J. e WILLIAMS, T.A., An Introduction to Management Science 9ª ed., South-Western College Publishing, 2000.`n4. PIZZOLATO, N. D. e GANDOLPHO, A. A. Técnicas de Otimização, LTC Editora, 2009.`n5. TAHA, H. A ., Pesquisa Operacional, 8ª ed., Pearson/Prentice Hall, 2008."
$ws.Cells.Item(22, 2).Value = $biblio
$ws.Cells.Item(22, 3).Value = $biblio
